# loading_percent.xlsx - case with 380 kV done
# Recalculated loading-percent results for rows 2-25 (columns C:H and M:N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "C2" = 4.924741707968425; "D2" = 4.522883239565808; "E2" = 9.427305210633097; "F2" = 32.78506309923144; "G2" = 47.70859717596404; "H2" = 16.19504843632827; "M2" = 23.63142822161213; "N2" = 16.70351180470752
    "C3" = 4.749893571430197; "D3" = 4.472065699592989; "E3" = 9.359589629406411; "F3" = 31.65327972735108; "G3" = 45.45912431382495; "H3" = 15.91597916326749; "M3" = 22.61520260742766; "N3" = 16.48521143422583
    "C4" = 4.641213922249149; "D4" = 4.442890097197288; "E4" = 9.319394675614801; "F4" = 30.95320123860136; "G4" = 44.04277661481686; "H4" = 15.74919324837453; "M4" = 21.9712839003748; "N4" = 16.35192556931985
    "C5" = 4.596670915058227; "D5" = 4.43151924372188; "E5" = 9.30338430182846; "F5" = 30.66713030263638; "G5" = 43.45772173858139; "H5" = 15.68246181679581; "M5" = 21.7042363708297; "N5" = 16.29786064558926
    "C6" = 4.58926154639036; "D6" = 4.429662662713857; "E6" = 9.300748726794488; "F6" = 30.61959681932184; "G6" = 43.36012835862014; "H6" = 15.67145821516002; "M6" = 21.65962457085954; "N6" = 16.28890010441759
    "C7" = 4.640614126307558; "D7" = 4.442734636190197; "E7" = 9.319177230116352; "F7" = 30.94934565298077; "G7" = 44.03491691006859; "H7" = 15.74828817020069; "M7" = 21.96770068658636; "N7" = 16.35119533708602
    "C8" = 4.864772999592326; "D8" = 4.504946507369328; "E8" = 9.403683073650958; "F8" = 32.39619242340822; "G8" = 46.94081396194422; "H8" = 16.0979315394264; "M8" = 23.28538300156355; "N8" = 16.62812252781465
    "C9" = 5.290774549612983; "D9" = 4.642563247068112; "E9" = 9.579393146678646; "F9" = 35.17179385879285; "G9" = 52.32433693321408; "H9" = 16.81582976004993; "M9" = 25.69706023233983; "N9" = 17.17446715302819
    "C10" = 5.591807441768802; "D10" = 4.752486609760331; "E10" = 9.713245784641771; "F10" = 37.1484550010579; "G10" = 56.0474601661676; "H10" = 17.35755308749641; "M10" = 27.34813918444793; "N10" = 17.5743461363387
    "C11" = 5.725510448221045; "D11" = 4.804236330418847; "E11" = 9.774884180233485; "F11" = 38.02962475540119; "G11" = 57.6843245950461; "H11" = 17.60596565598748; "M11" = 28.07045869709336; "N11" = 17.75520627501886
    "C12" = 5.775628634910023; "D12" = 4.824068045507966; "E12" = 9.798308493813479; "F12" = 38.36038907455616; "G12" = 58.29556836148655; "H12" = 17.70022560498526; "M12" = 28.33967413577568; "N12" = 17.82348433417746
    "C13" = 5.764858250103676; "D13" = 4.819786724619754; "E13" = 9.793260335851983; "F13" = 38.2892874876122; "G13" = 58.16431474421839; "H13" = 17.67991797340099; "M13" = 28.2818880118328; "N13" = 17.80878968271274
    "C14" = 5.729644214224194; "D14" = 4.805863275115473; "E14" = 9.776809784634859; "F14" = 38.05689708823243; "G14" = 57.73478658998981; "H14" = 17.61371716407389; "M14" = 28.09269452138248; "N14" = 17.76082808783361
    "C15" = 5.708006606676316; "D15" = 4.797364918068394; "E15" = 9.766743458963473; "F15" = 37.91416260011648; "G15" = 57.47055645768854; "H15" = 17.57318947991366; "M15" = 27.97624187221477; "N15" = 17.73142118689663
    "C16" = 5.583000360205674; "D16" = 4.749138426709154; "E16" = 9.709230537783894; "F16" = 37.09047726598624; "G16" = 55.93930547806222; "H16" = 17.34135054994458; "M16" = 27.30033948509873; "N16" = 17.56250032929849
    "C17" = 5.505448690892701; "D17" = 4.719988970218219; "E17" = 9.674122459384979; "F17" = 36.58031432449487; "G17" = 54.98505180844785; "H17" = 17.1995647819276; "M17" = 26.87819436706241; "N17" = 17.45856199091901
    "C18" = 5.460540618626617; "D18" = 4.703388204986362; "E18" = 9.654001097752769; "F18" = 36.28520254444883; "G18" = 54.43086455453934; "H18" = 17.11820513962563; "M18" = 26.63268873500008; "N18" = 17.39868430231884
    "C19" = 5.445285092099325; "D19" = 4.697796319941166; "E19" = 9.647201440098534; "F19" = 36.18500479032872; "G19" = 54.24232599305835; "H19" = 17.09069398718681; "M19" = 26.54910700143491; "N19" = 17.3783962476201
    "C20" = 5.513735851522619; "D20" = 4.723074991059339; "E20" = 9.677852512375663; "F20" = 36.63479833306472; "G20" = 55.08718831164365; "H20" = 17.21463895472085; "M20" = 26.92341312083276; "N20" = 17.46963668276253
    "C21" = 5.740001675242061; "D21" = 4.809946673569562; "E21" = 9.781639646301358; "F21" = 38.12523730338219; "G21" = 57.86118594507962; "H21" = 17.6331574772443; "M21" = 28.14838349949619; "N21" = 17.77492172543205
    "C22" = 5.88487606510522; "D22" = 4.868085793636958; "E22" = 9.849944229646907; "F22" = 39.08220716125638; "G22" = 59.6238679091458; "H22" = 17.90775875182701; "M22" = 28.92378073522909; "N22" = 17.97319455027925
    "C23" = 5.807842625392327; "D23" = 4.836936400139169; "E23" = 9.813453467672071; "F23" = 38.57311802848824; "G23" = 58.68781708321501; "H23" = 17.76113049063084; "M23" = 28.51229246160891; "N23" = 17.86750589626934
    "C24" = 5.509990227385981; "D24" = 4.721679308819794; "E24" = 9.676165958145976; "F24" = 36.6101717327528; "G24" = 55.04102972702592; "H24" = 17.20782343166588; "M24" = 26.90297845741125; "N24" = 17.46463019054724
    "C25" = 5.177387271956397; "D25" = 4.603730516096646; "E25" = 9.530925801342786; "F25" = 34.4302559039676; "G25" = 50.90621974918423; "H25" = 16.61872161205481; "M25" = 25.06486081643276; "N25" = 17.02668227744613
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
